# Sari_R3_9.xlsx template update
# - Fill in the "Semana epidemiologica" (week number) column B with 1..53
# - Fix the C6 header cell style so it wraps text like its neighbours
# - Move the active selection from B3 to the header row B1:O1
# - Give the bar chart's value axis a title ("Numero de casos SARI")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data: populate the epidemiological week numbers (B7:B59) ---
for ($i = 0; $i -lt 53; $i++) {
    $ws.Cells.Item(7 + $i, 2).Value = $i + 1
}

# --- Header cell C6: turn on wrap text to match the other header cells ---
$ws.Range("C6").WrapText = $true

# --- Selection moves to the title row ---
$ws.Range("B1:O1").Select()

# --- Chart: add a title to the value (Y) axis ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.AxisTitle.Text = "Número de casos SARI"
$valAx.HasTitle = $true
